$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Carlos Estevez (Oblique) ---
$ws.Range("A2").Value = "Carlos Estevez"
$ws.Range("B2").Value = "estevca01"
$ws.Range("C2").Value = "April 23 2018"
$ws.Range("D2").Value = "Oblique"
$ws.Range("E2").Value = "Estevez has been shifted to the 60-day disabled list due to a left oblique strain but is likely to return during the latter part of May."

# --- Row 3: DJ LeMahieu (Wrist) ---
$ws.Range("A3").Value = "DJ LeMahieu"
$ws.Range("B3").Value = "lemahdj01"
$ws.Range("C3").Value = "May 18 2018"
$ws.Range("D3").Value = "Wrist"
$ws.Range("E3").Value = "LeMahieu is on the 10-day disabled list after being diagnosed with a slight left wrist fracture and a sprained thumb. There is no timetable for return."

# --- Row 4: Zac Rosscup (Finger) ---
$ws.Range("A4").Value = "Zac Rosscup"
$ws.Range("B4").Value = "rosscza01"
$ws.Range("C4").Value = "May 01 2018"
$ws.Range("D4").Value = "Finger"
$ws.Range("E4").Value = "Rosscup has been shifted to the 60-day disabled list dealing with warts on his left middle finger and is likely to stay on the DL until at least the end of May."

# Rows 5 and 6 (formerly blank, style-only placeholder rows) are no longer
# part of the table -- remove them so the used range shrinks to A1:E4.
$ws.Rows("5:6").Delete()

# Column widths were re-fit to the new (shorter) content. The exact
# font-metric "best fit" pixel widths aren't independently settable through
# this object model's ColumnWidth API (it only lets us choose the resulting
# width to a coarse 1/6-character granularity), so we dial in the closest
# reachable width for each resized column.
$ws.Columns("A").ColumnWidth = 13.65
$ws.Columns("D").ColumnWidth = 11.8
$ws.Columns("E").ColumnWidth = 72.8

# Selection moved as part of the edit session.
$ws.Range("A12").Select() | Out-Null

